# Comentario antes de clase
# Replace the three intro paragraphs (hyperlink, blank line, "Creo en archivo...")
# and the two closing paragraphs with the new class-prep notes, preserving the
# trailing sectPr untouched. We rebuild the whole body range via InsertXML so
# that the paragraph count collapses from 5 to 4 in one shot.

$d = $word.ActiveDocument

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Preguntar si código de muestra de videos est</w:t></w:r><w:r><w:t>á</w:t></w:r><w:r><w:t xml:space="preserve"> bien hecho</w:t></w:r><w:r><w:t xml:space="preserve"> y como implementarlo</w:t></w:r><w:r><w:t xml:space="preserve">. Además de ver como se quedaría el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>html</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de muestra de videos.</w:t></w:r></w:p><w:p><w:r><w:t>Preguntar s la manera de introducir datos fijos en el formulario</w:t></w:r><w:r><w:t xml:space="preserve">, es decir, el correo y la contraseña que nos ofrecen. (He estado investigando y la única opción que había encontrado es utilizando </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>la opciones</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> de validaciones de formulario)</w:t></w:r></w:p><w:p><w:r><w:t>Por ultimo ver si falta algo mas aunque esta casi todo terminado.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml) | Out-Null
